# river update May 2024
# Append new water-quality sample rows (160-172) for
# "Tokiahuru at Karioi Domain Road" to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$siteName = "Tokiahuru at Karioi Domain Road"

# Each entry: row, parameter name, date-serial (Excel date serial, matching the
# existing "date time" column), Value (kept as text - this sheet stores the
# Value column as text so things like "0.405"/"70.000" keep their formatting),
# and Unit (blank "" for the macroinvertebrate-index rows).
$rows = @(
    @{ Row = 160; Param = "Chlorophyll A";                                         Date = 44761.00046296296; Value = "0.405";  Unit = "mg/m2" },
    @{ Row = 161; Param = "Chlorophyll A";                                         Date = 44789.00046296296; Value = "6.500";  Unit = "mg/m2" },
    @{ Row = 162; Param = "Chlorophyll A";                                         Date = 44852.00046296296; Value = "4.650";  Unit = "mg/m2" },
    @{ Row = 163; Param = "Chlorophyll A";                                         Date = 44887.00046296296; Value = "4.350";  Unit = "mg/m2" },
    @{ Row = 164; Param = "Chlorophyll A";                                         Date = 44914.00046296296; Value = "2.450";  Unit = "mg/m2" },
    @{ Row = 165; Param = "Chlorophyll A";                                         Date = 44950.00046296296; Value = "5.500";  Unit = "mg/m2" },
    @{ Row = 166; Param = "Chlorophyll A";                                         Date = 44978.00046296296; Value = "1.550";  Unit = "mg/m2" },
    @{ Row = 167; Param = "ASPM (Macroinvertebrate Average Score Per Metric)";     Date = 45020;             Value = "0.389";  Unit = "" },
    @{ Row = 168; Param = "MCI (Macroinvertebrate Community Index)";               Date = 45020;             Value = "104.44"; Unit = "" },
    @{ Row = 169; Param = "QMCI (Quantitative Macroinvertebrate Community Index)"; Date = 45020;             Value = "5.180";  Unit = "" },
    @{ Row = 170; Param = "Chlorophyll A";                                         Date = 45042.00046296296; Value = "70.000"; Unit = "mg/m2" },
    @{ Row = 171; Param = "Chlorophyll A";                                         Date = 45062.00046296296; Value = "6.000";  Unit = "mg/m2" },
    @{ Row = 172; Param = "Chlorophyll A";                                         Date = 45097.00046296296; Value = "29.500"; Unit = "mg/m2" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # A: site name (plain text)
    $ws.Cells.Item($row, 1).Value = $siteName

    # B: parameter name (plain text)
    $ws.Cells.Item($row, 2).Value = $r.Param

    # C: date time - numeric date serial, displayed the same way as the rest
    # of the column
    $cDate = $ws.Cells.Item($row, 3)
    $cDate.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cDate.Value = $r.Date

    # D: Value - stored as text (leading "'" forces text so "0.405"/"70.000"
    # keep their exact digits/trailing zeros instead of becoming numbers)
    $ws.Cells.Item($row, 4).Value = "'" + $r.Value

    # E: Project (always blank text in this sheet)
    $ws.Cells.Item($row, 5).Value = "'"

    # F: Method (always blank text in this sheet)
    $ws.Cells.Item($row, 6).Value = "'"

    # G: Unit (blank text for index rows, otherwise a plain unit string)
    if ($r.Unit -eq "") {
        $ws.Cells.Item($row, 7).Value = "'"
    } else {
        $ws.Cells.Item($row, 7).Value = $r.Unit
    }

    # H: Quality code (numeric)
    $ws.Cells.Item($row, 8).Value = 200

    # I: pH (always blank text in this sheet)
    $ws.Cells.Item($row, 9).Value = "'"
}
